$d = $word.ActiveDocument

# 1. recuite -> recuit
$d.Content.Find.Execute("recuite", $false, $false, $false, $false, $false, $true, 1, $false, "recuit", 2) | Out-Null

# 2. nest pas pur car ony mesle -> nest pas pur car on y mesle
$d.Content.Find.Execute("ony mesle", $false, $false, $false, $false, $false, $true, 1, $false, "on y mesle", 2) | Out-Null

# 3. Mays les prendre pur & -> Mays les prendre purs &
$d.Content.Find.Execute("prendre pur &", $true, $false, $false, $false, $false, $true, 1, $false, "prendre purs &", 2) | Out-Null

# 4. estant en -> estans en  (inserted as a separate run "s" replacing "t")
$d.Content.Find.Execute("estant en", $false, $false, $false, $false, $false, $true, 1, $false, "estans en", 2) | Out-Null

# 5. ou subtillie avecq l -> ou subtilies avecq l
$d.Content.Find.Execute("subtillie", $false, $false, $false, $false, $false, $true, 1, $false, "subtilies", 2) | Out-Null
